# Updates the crypto price/volume table (columns D & E) for the Tue Oct 29
# 2024 GitHub Actions refresh, plus the row-19/20 swap (Chainlink <-> BitcoinCash).
# Cells whose new text is a plain decimal number (e.g. "605.46") are forced to
# stay text (matching the original inlineStr cells) by briefly switching the
# cell to a Text number format, assigning the value, then restoring the
# "Normal" style so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.666.06"
$ws.Range("E2").Value = "  +4.11%  "
$ws.Range("D3").Value = "2.634.26"
$ws.Range("E3").Value = "  +2.70%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.51%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +1.72%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.175"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +8.87%  "
$ws.Range("D10").Value = "2.633.32"
$ws.Range("E10").Value = "  +2.68%  "
$ws.Range("E11").Value = "  +1.32%  "
$ws.Range("E12").Value = "  +3.52%  "
$ws.Range("E13").Value = "  +0.47%  "
$ws.Range("E14").Value = "  +4.70%  "
$ws.Range("D15").Value = "3.108.02"
$ws.Range("D16").Value = "72.603.41"
$ws.Range("E16").Value = "  +4.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.84"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.08%  "
$ws.Range("D18").Value = "2.622.24"
$ws.Range("E18").Value = "  +1.85%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.48%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "385.74"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.49%  "
$ws.Range("E21").Value = "  +3.13%  "
$ws.Range("E22").Value = "  +1.57%  "
$ws.Range("E23").Value = "  +16.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.83%  "
$ws.Range("E25").Value = "  +2.79%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.75%  "
$ws.Range("D28").Value = "2.730.21"
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").Value = "0.0₃0960"
$ws.Range("E30").Value = "  +4.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "521.44"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.97%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.12"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.35%  "
$ws.Range("E33").Value = "  +4.14%  "
$ws.Range("E34").Value = "  +1.71%  "
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "164.73"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.42"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.11%  "
$ws.Range("E38").Value = "  +4.33%  "
$ws.Range("E39").Value = "  +1.04%  "
$ws.Range("E40").Value = "  -5.81%  "
$ws.Range("E41").Value = "  +5.55%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.17"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.93%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("E44").Value = "  +4.70%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.335"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.78%  "
$ws.Range("E46").Value = "  +0.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "150.99"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.81%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.71"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.14%  "
$ws.Range("E49").Value = "  +4.49%  "
$ws.Range("E50").Value = "  +5.14%  "
$ws.Range("E51").Value = "  +2.67%  "
